# Enhance Azure DevOps integration and improve chatbot functionality
#
# Sheet "This Friday" (sheet1):
#   - Remove the "TEST_03: Coffee Cake" work item (ID 986)
#   - Move TEST_01/TEST_02/TEST_04/TEST_05 items to "In Progress"
#   - Assign TEST_05: Strawberry (ID 988) to Priththiha Nemikumar
#
# Sheet "Next Friday" (sheet2):
#   - Remove the "QA | Sure, whatever" work item (ID 945)
#
# Sheet "Friday After Next" (sheet3):
#   - Replace the "Arian Fooladray" display name with his e-mail address
#   - Move "Bug | Minor | Performance Degradation" (ID 940) to "In Progress"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "This Friday"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("This Friday")

# Row 4 holds ID 986 ("TEST_03: Coffee Cake") - remove it entirely, shifting
# the rows below it up.
$ws1.Rows.Item(4).Delete()

# After the deletion the remaining TEST_* rows occupy rows 2, 3, 4 and 5.
# Row 2 = ID 988 "TEST_05: Strawberry"
$ws1.Range("C2").Value = "In Progress"
$ws1.Range("D2").Value = "Priththiha Nemikumar"

# Row 3 = ID 987 "TEST_04: Testing Sandwiches"
$ws1.Range("C3").Value = "In Progress"

# Row 4 = ID 985 "TEST_02: Something"
$ws1.Range("C4").Value = "In Progress"

# Row 5 = ID 984 "TEST_01: Introduction"
$ws1.Range("C5").Value = "In Progress"

# ---------------------------------------------------------------------
# Sheet 2: "Next Friday"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Next Friday")

# Row 2 holds ID 945 ("QA | Sure, whatever") - remove it entirely, shifting
# the rows below it up.
$ws2.Rows.Item(2).Delete()

# ---------------------------------------------------------------------
# Sheet 3: "Friday After Next"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Friday After Next")

# Swap the "Arian Fooladray" display name for his e-mail address wherever
# it appears as the assignee (rows 3-8).
for ($r = 3; $r -le 8; $r++) {
    $cell = $ws3.Cells.Item($r, 4)
    if ($cell.Value2 -eq "Arian Fooladray") {
        $cell.Value = "afooladray@fgfbrands.com"
    }
}

# ID 940 ("Bug | Minor | Performance Degradation") moves to "In Progress".
$ws3.Range("C2").Value = "In Progress"
